$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "66.095.32"
$ws.Range("E2").Value = "  +1.65%  "

# Row 3
$ws.Range("D3").Value = "2.694.26"
$ws.Range("E3").Value = "  +2.65%  "

# Row 4
$ws.Range("E4").Value = "  +0.03%  "

# Row 5
$ws.Range("D5").Value = "'616.38"
$ws.Range("E5").Value = "  +2.45%  "

# Row 6
$ws.Range("D6").Value = "'158.44"
$ws.Range("E6").Value = "  +2.45%  "

# Row 7
$ws.Range("E7").Value = "  +0.03%  "

# Row 8
$ws.Range("D8").Value = "'0.591"
$ws.Range("E8").Value = "  +0.88%  "

# Row 9
$ws.Range("D9").Value = "'0.125"
$ws.Range("E9").Value = "  +6.33%  "

# Row 10
$ws.Range("D10").Value = "'6.04"
$ws.Range("E10").Value = "  +4.89%  "

# Row 11
$ws.Range("D11").Value = "'0.404"
$ws.Range("E11").Value = "  -1.01%  "

# Row 12
$ws.Range("E12").Value = "  +0.82%  "

# Row 13
$ws.Range("D13").Value = "'0.0000210"
$ws.Range("E13").Value = "  +10.98%  "

# Row 14
$ws.Range("D14").Value = "'30.20"
$ws.Range("E14").Value = "  +3.48%  "

# Row 15
$ws.Range("D15").Value = "3.180.95"
$ws.Range("E15").Value = "  +2.74%  "

# Row 16
$ws.Range("D16").Value = "65.949.70"
$ws.Range("E16").Value = "  +1.55%  "

# Row 17
$ws.Range("D17").Value = "2.696.07"
$ws.Range("E17").Value = "  +2.39%  "

# Row 18
$ws.Range("D18").Value = "'12.72"
$ws.Range("E18").Value = "  +1.39%  "

# Row 19
$ws.Range("D19").Value = "'4.90"
$ws.Range("E19").Value = "  +0.21%  "

# Row 20
$ws.Range("D20").Value = "'7.80"
$ws.Range("E20").Value = "  +6.61%  "

# Row 21
$ws.Range("D21").Value = "'358.56"
$ws.Range("E21").Value = "  -0.03%  "

# Row 22
$ws.Range("B22").Value = "Dai"
$ws.Range("C22").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D22").Value = "'0.999"
$ws.Range("E22").Value = "  -0.07%  "

# Row 23
$ws.Range("B23").Value = "Litecoin"
$ws.Range("C23").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D23").Value = "'71.18"
$ws.Range("E23").Value = "  +3.59%  "

# Row 24
$ws.Range("D24").Value = "'0.0000113"
$ws.Range("E24").Value = "  +19.06%  "

# Row 25
$ws.Range("D25").Value = "'9.91"
$ws.Range("E25").Value = "  +5.81%  "

# Row 26
$ws.Range("D26").Value = "'1.64"
$ws.Range("E26").Value = "  -0.67%  "

# Row 27
$ws.Range("D27").Value = "'1.67"
$ws.Range("E27").Value = "  +2.58%  "

# Row 28
$ws.Range("D28").Value = "'0.172"
$ws.Range("E28").Value = "  +4.27%  "

# Row 29
$ws.Range("D29").Value = "'8.30"
$ws.Range("E29").Value = "  +1.66%  "

# Row 30
$ws.Range("E30").Value = "  +5.41%  "

# Row 31
$ws.Range("B31").Value = "Binance-PegBSC-USD"
$ws.Range("C31").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D31").Value = "'0.999"
$ws.Range("E31").Value = "  -0.01%  "

# Row 32
$ws.Range("B32").Value = "Bittensor"
$ws.Range("C32").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D32").Value = "'535.54"
$ws.Range("E32").Value = "  +2.42%  "

# Row 33
$ws.Range("D33").Value = "'1.81"
$ws.Range("E33").Value = "  +0.26%  "

# Row 34
$ws.Range("D34").Value = "'6.70"
$ws.Range("E34").Value = "  +5.83%  "

# Row 35
$ws.Range("D35").Value = "'5.46"
$ws.Range("E35").Value = "  -0.34%  "

# Row 36
$ws.Range("D36").Value = "'0.436"
$ws.Range("E36").Value = "  +2.40%  "

# Row 37
$ws.Range("D37").Value = "'20.82"
$ws.Range("E37").Value = "  +2.25%  "

# Row 38
$ws.Range("D38").Value = "'164.80"
$ws.Range("E38").Value = "  +1.68%  "

# Row 39
$ws.Range("D39").Value = "'2.00"
$ws.Range("E39").Value = "  -0.44%  "

# Row 40
$ws.Range("D40").Value = "'1.00"
$ws.Range("E40").Value = "  +0.06%  "

# Row 41
$ws.Range("B41").Value = "USDe"
$ws.Range("C41").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D41").Value = "'0.999"
$ws.Range("E41").Value = "  -0.03%  "

# Row 42
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").Value = "'169.25"
$ws.Range("E42").Value = "  +2.92%  "

# Row 43
$ws.Range("D43").Value = "'42.53"
$ws.Range("E43").Value = "  +0.85%  "

# Row 44
$ws.Range("D44").Value = "'4.19"
$ws.Range("E44").Value = "  +1.13%  "

# Row 45
$ws.Range("D45").Value = "'0.0626"
$ws.Range("E45").Value = "  +2.30%  "

# Row 46
$ws.Range("D46").Value = "'2.33"
$ws.Range("E46").Value = "  +5.31%  "

# Row 47
$ws.Range("D47").Value = "'23.73"
$ws.Range("E47").Value = "  +2.40%  "

# Row 48
$ws.Range("E48").Value = "  +2.11%  "

# Row 49
$ws.Range("D49").Value = "'0.658"
$ws.Range("E49").Value = "  +0.91%  "

# Row 50
$ws.Range("D50").Value = "'21.00"
$ws.Range("E50").Value = "  +7.71%  "

# Row 51
$ws.Range("D51").Value = "'0.0995"
$ws.Range("E51").Value = "  +1.58%  "
